# CasosColombia.xlsx - "Add files via upload" update
#
# This mirrors a re-upload of the tracked daily COVID-19 case file:
#   - A handful of historical cells that had been mistakenly logged with a
#     numeric placeholder got normalized back to the "NaN" text marker that
#     the rest of the sheet uses for missing data (and one cell that had
#     drifted to "NaN" got restored to its real numeric reading).
#   - A new day's row (row 194, 2020-09-14) was appended at the bottom of
#     the data table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Fix up a handful of mis-encoded historical cells -----------------
# These columns use the shared "NaN" text marker elsewhere in the sheet for
# missing/unavailable readings; these particular cells had a stray number
# instead, so restore the "NaN" marker.
$ws.Range("CZ18").Value = "NaN"
$ws.Range("AD25").Value = "NaN"
$ws.Range("CQ34").Value = "NaN"
$ws.Range("AD66").Value = "NaN"
$ws.Range("AD67").Value = "NaN"
$ws.Range("AD68").Value = "NaN"
$ws.Range("AD192").Value = "NaN"

# ... and this one had incorrectly been left as "NaN" instead of the real
# reported count.
$ws.Range("CR31").Value = 1

# --- Append the new day's data row (row 194) ---------------------------
$row194 = @(
    44088, 721892, 2725, 96560, 65843, 242007, 27277, 5226, 4078, 7319,
    7348, 15758, 3865, 22132, 28588, 6676, 7706, 13871, 11850, 16194,
    13682, 3380, 2091, 8316, 24927, 13242, 9651, 53927, 1565, 591,
    589, 461, 349, 268, 473, 2005, 4405, 37108, 8401, 2505,
    42089, 1046, 21938, 1495, 9411, 1610, 1589, 6430, 1782, 952,
    2484, 2653, 55857, 13347, 4502, 8727, 5547, 281, 1436, 2654,
    737, 2123, 9223, 9205, 9832, 14092, 1932, 864, 11676, 9399,
    10973, 2145, 1855, 4550, 4216, 1493, 5349, 3075, 1767, 852,
    2615, 2130, 1681, 1308, 5928, 1853, 1325, 1597, 1945, 1918,
    2259, 1418, 1182, 1163, 783, 3194, 1301, 886, 924, 1630,
    1455, 719, 826, 1155, 1433, 1249, 1355, 1067, 327, 352,
    768, 688, 457, 535, 364, 652, 737, 519, 485, 372,
    518, 129425, 306749, 14758, 131924, 81334, 38658, 11048
)

$newRow = 194
for ($i = 0; $i -lt $row194.Length; $i++) {
    $ws.Cells.Item($newRow, $i + 1).Value = $row194[$i]
}

# Move the selection down to the new last cell, same as Excel would leave
# it after typing in the last value of the new row.
$ws.Range("DX194").Select()
